# Generate Report for Handback
#
# This script mirrors a "handback" run: the localized files for zh-cn/de-de
# are now in sync with en-US, so the status text is updated, the
# "Latest Target File"/"Latest Handback File"/"Latest Handback DateTime"
# columns get populated for row 2 and row 3 on both locale sheets, and a
# hyperlink to the target markdown file ("a.md") is added for the
# "Latest Target File" cell.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status columns for zh-cn (E) and de-de (F), rows 2-3 ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status column (C), rows 2-3 ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# --- de-de sheet: Status column (C), rows 2-3 ---
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8247c9ae795d26dfd56c36dc98d0d0044dc51c2f/e2e/a.md"

# --- zh-cn sheet: Latest Target File (I), Latest Handback File (J),
#     Latest Handback DateTime (K) for rows 2-3 ---
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $targetUrl, "", "", "a.md")
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-27 22:37:44"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $targetUrl, "", "", "a.md")
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-27 22:37:44"

# --- de-de sheet: Latest Target File (I), Latest Handback File (J),
#     Latest Handback DateTime (K) for rows 2-3 ---
$dede.Hyperlinks.Add($dede.Range("I2"), $targetUrl, "", "", "a.md")
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-27 22:37:51"

$dede.Hyperlinks.Add($dede.Range("I3"), $targetUrl, "", "", "a.md")
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-08-27 22:37:51"

# --- Widen the columns that now hold the longer status / file-name text,
#     matching Excel's autofit-on-edit behaviour. ---
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(10).ColumnWidth = 40

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(10).ColumnWidth = 40
